$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 397.45715
$ws.Range("I28").Value = 373.61905
$ws.Range("J28").Value = 433.2143
$ws.Range("K28").Value = 373.61905
$ws.Range("L28").Value = 433.2143
$ws.Range("M28").Value = 111.38095
$ws.Range("N28").Value = -1403.2143
# Row 62
$ws.Range("H62").Value = 62501780
$ws.Range("I62").Value = 76924790
$ws.Range("J62").Value = 2053.3333
$ws.Range("K62").Value = 76924790
$ws.Range("L62").Value = 2053.3333
$ws.Range("M62").Value = -76924166
$ws.Range("N62").Value = -3301.3333
# Row 65
$ws.Range("H65").Value = 62501780
$ws.Range("I65").Value = 76924790
$ws.Range("J65").Value = 2053.3333
$ws.Range("K65").Value = 384623950
$ws.Range("L65").Value = 10266.6665
$ws.Range("M65").Value = -384620830
$ws.Range("N65").Value = -16506.6665
# Row 76
$ws.Range("H76").Value = 74510.39
$ws.Range("I76").Value = 95289.21000000001
$ws.Range("J76").Value = 3268.7144
$ws.Range("K76").Value = 95289.21000000001
$ws.Range("L76").Value = 3268.7144
$ws.Range("M76").Value = -94974.21000000001
$ws.Range("N76").Value = -3898.7144
# Row 79
$ws.Range("H79").Value = 74510.39
$ws.Range("I79").Value = 95289.21000000001
$ws.Range("J79").Value = 3268.7144
$ws.Range("K79").Value = 95289.21000000001
$ws.Range("L79").Value = 3268.7144
$ws.Range("M79").Value = -94197.21000000001
$ws.Range("N79").Value = -5452.7144
# Row 86
$ws.Range("H86").Value = 7684210.5
$ws.Range("I86").Value = 18779248
$ws.Range("J86").Value = 3030.7693
$ws.Range("K86").Value = 18779248
$ws.Range("L86").Value = 3030.7693
$ws.Range("M86").Value = -18778125
$ws.Range("N86").Value = -5276.7693
# Row 89
$ws.Range("H89").Value = 7684210.5
$ws.Range("I89").Value = 18779248
$ws.Range("J89").Value = 3030.7693
$ws.Range("K89").Value = 93896240
$ws.Range("L89").Value = 15153.8465
$ws.Range("M89").Value = -93890624
$ws.Range("N89").Value = -26385.8465
# Row 92
$ws.Range("H92").Value = 5208955
$ws.Range("I92").Value = 7576191
$ws.Range("J92").Value = 1036.4
$ws.Range("K92").Value = 7576191
$ws.Range("L92").Value = 1036.4
$ws.Range("M92").Value = -7574943
$ws.Range("N92").Value = -3532.4
# Row 98
$ws.Range("H98").Value = 942.94116
$ws.Range("I98").Value = 942.94116
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 942.94116
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 555.05884
$ws.Range("N98").ClearContents()
# Row 122
$ws.Range("H122").Value = 942.94116
$ws.Range("I122").Value = 942.94116
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2828.82348
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -378.82348
$ws.Range("N122").ClearContents()
# Row 137
$ws.Range("H137").Value = 40001324
$ws.Range("I137").Value = 27028022
$ws.Range("J137").Value = 76925340
$ws.Range("K137").Value = 81084066
$ws.Range("L137").Value = 230776020
$ws.Range("M137").Value = -81081516
$ws.Range("N137").Value = -230781120
# Row 138
$ws.Range("H138").Value = 1945.6364
$ws.Range("I138").Value = 682.1515000000001
$ws.Range("J138").Value = 2893.25
$ws.Range("K138").Value = 2046.4545
$ws.Range("L138").Value = 8679.75
$ws.Range("M138").Value = 3093.5455
$ws.Range("N138").Value = -18959.75

$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 18000
$ws.Range("J23").Value = 18000
$ws.Range("L23").Value = 18000
$ws.Range("N23").Value = -18518

$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 288.84616
$ws.Range("I19").Value = 135.55556
$ws.Range("J19").Value = 633.75
$ws.Range("K19").Value = 135.55556
$ws.Range("L19").Value = 633.75
$ws.Range("M19").Value = 34.44443999999999
$ws.Range("N19").Value = -973.75
# Row 24
$ws.Range("H24").Value = 288.84616
$ws.Range("I24").Value = 135.55556
$ws.Range("J24").Value = 633.75
$ws.Range("K24").Value = 135.55556
$ws.Range("L24").Value = 633.75
$ws.Range("M24").Value = 34.44443999999999
$ws.Range("N24").Value = -973.75
# Row 58
$ws.Range("H58").Value = 917.9636
$ws.Range("I58").Value = 476.3243
$ws.Range("J58").Value = 1825.7778
$ws.Range("K58").Value = 476.3243
$ws.Range("L58").Value = 1825.7778
$ws.Range("M58").Value = -273.3243
$ws.Range("N58").Value = -2231.7778
# Row 136
$ws.Range("H136").Value = 917.9636
$ws.Range("I136").Value = 476.3243
$ws.Range("J136").Value = 1825.7778
$ws.Range("K136").Value = 1428.9729
$ws.Range("L136").Value = 5477.3334
$ws.Range("M136").Value = 1121.0271
$ws.Range("N136").Value = -10577.3334

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 717.12
$ws.Range("I5").Value = 636
$ws.Range("J5").Value = 1650
$ws.Range("K5").Value = 1908
$ws.Range("L5").Value = 4950
$ws.Range("M5").Value = -1796
$ws.Range("N5").Value = -5174
# Row 9
$ws.Range("H9").Value = 185500.33
$ws.Range("I9").Value = 220600.4
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 661801.2
$ws.Range("L9").Value = 30000
$ws.Range("M9").Value = -661577.2
$ws.Range("N9").Value = -30448
# Row 122
$ws.Range("H122").Value = 654.1539
$ws.Range("I122").Value = 400
$ws.Range("J122").Value = 1060.8
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 9547.199999999999
$ws.Range("M122").Value = -1150
$ws.Range("N122").Value = -14447.2
# Row 127
$ws.Range("H127").Value = 1337
$ws.Range("J127").Value = 1337
$ws.Range("L127").Value = 4011
$ws.Range("N127").Value = -13931
# Row 135
$ws.Range("H135").Value = 717.12
$ws.Range("I135").Value = 636
$ws.Range("J135").Value = 1650
$ws.Range("K135").Value = 5724
$ws.Range("L135").Value = 14850
$ws.Range("M135").Value = -3189
$ws.Range("N135").Value = -19920

$ws = $wb.Worksheets.Item("GSM")
# Row 103
$ws.Range("H103").Value = 25360
$ws.Range("J103").Value = 25360
$ws.Range("L103").Value = 25360
$ws.Range("N103").Value = -27704

$ws = $wb.Worksheets.Item("LTW")
# Row 23
$ws.Range("H23").Value = 21265.8
$ws.Range("I23").Value = 21265.8
$ws.Range("K23").Value = 21265.8
$ws.Range("M23").Value = -21035.8
# Row 68
$ws.Range("H68").Value = 16113556
$ws.Range("I68").Value = 67667770
$ws.Range("J68").Value = 2864.5
$ws.Range("K68").Value = 67667770
$ws.Range("L68").Value = 2864.5
$ws.Range("M68").Value = -67667021
$ws.Range("N68").Value = -4362.5
# Row 71
$ws.Range("H71").Value = 16113556
$ws.Range("I71").Value = 67667770
$ws.Range("J71").Value = 2864.5
$ws.Range("K71").Value = 338338850
$ws.Range("L71").Value = 14322.5
$ws.Range("M71").Value = -338335106
$ws.Range("N71").Value = -21810.5
# Row 107
$ws.Range("H107").Value = 25000
$ws.Range("I107").Value = 25000
$ws.Range("K107").Value = 25000
$ws.Range("M107").Value = -23080

$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("H47").Value = 14031
$ws.Range("I47").Value = 10062
$ws.Range("J47").Value = 18000
$ws.Range("K47").Value = 10062
$ws.Range("L47").Value = 18000
$ws.Range("M47").Value = -9490
$ws.Range("N47").Value = -19144
# Row 107
$ws.Range("H107").Value = 672.6875
$ws.Range("I107").Value = 407.5
$ws.Range("J107").Value = 1468.25
$ws.Range("K107").Value = 1222.5
$ws.Range("L107").Value = 4404.75
$ws.Range("M107").Value = 697.5
$ws.Range("N107").Value = -8244.75
# Row 136
$ws.Range("H136").Value = 765.05884
$ws.Range("I136").Value = 408.01923
$ws.Range("J136").Value = 1327.6666
$ws.Range("K136").Value = 1224.05769
$ws.Range("L136").Value = 3982.9998
$ws.Range("M136").Value = 1325.94231
$ws.Range("N136").Value = -9082.9998
